$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")
Write-Host $ws.Name
